$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.504.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.86%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.645.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.99%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.003"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "302.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3839"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.77%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3595"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.40%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.97"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08165"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.227"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.64%  "

$ws.Range("E14").Value = "  +0.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.430"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.23%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001217"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.640.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.79%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "97.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.92%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07012"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.16%  "

$ws.Range("E20").Value = "  +2.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.55%  "

$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.506.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.479"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.87%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.027"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.232"
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.833.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.055"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.95%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.250"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.76%  "

$ws.Range("E35").Value = "  -2.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02786"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2496"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08773"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.048"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06968"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.26%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6965"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.333"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.81%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6482"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.61%  "

$ws.Range("E46").Value = "  +0.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.288"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.40%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.955"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07871"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.35%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.174"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.30%  "
